# Update "想去人数" (F) and "最低票价" (G) figures on the two sheets that
# carry the full exhibition listing: "展览" (sheet1) and "全部类型" (sheet4).
# The "全部类型" sheet has one extra row (a 演出/concert entry) inserted at
# row 3, so rows from 3 downward on "全部类型" are offset by +1 relative to
# the matching rows on "展览" (row 2 lines up 1:1 on both sheets).

$wb = $excel.ActiveWorkbook

# F-column (想去人数) updates, keyed by row number on the "展览" sheet.
$updates = @{
    2  = 1922
    6  = 2835
    7  = 197
    10 = 1588
    11 = 567
    12 = 49
    13 = 344
    15 = 26
    16 = 182
    23 = 22
    24 = 250
    25 = 52
    26 = 75
    27 = 1809
    29 = 432
    30 = 100
    33 = 319
    34 = 463
}

# G-column (最低票价) updates, keyed by the same "展览"-sheet row numbers.
$gUpdates = @{
    25 = 55
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # "全部类型" has an extra row inserted AT index 3 (a 演出/concert entry
    # that isn't part of "展览"), so row 2 lines up 1:1 between the two
    # sheets, but every "展览" row from 3 onward maps to row+1 on "全部类型".
    foreach ($row in $updates.Keys) {
        $targetRow = $row
        if (($sheetName -eq "全部类型") -and ($row -ge 3)) {
            $targetRow = $row + 1
        }
        $ws.Cells.Item($targetRow, 6).Value = $updates[$row]
    }

    foreach ($row in $gUpdates.Keys) {
        $targetRow = $row
        if (($sheetName -eq "全部类型") -and ($row -ge 3)) {
            $targetRow = $row + 1
        }
        $ws.Cells.Item($targetRow, 7).Value = $gUpdates[$row]
    }
}
